$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.530.50"
$ws.Range("E2").Value = "  +1.62%  "
$ws.Range("D3").Value = "2.585.67"
$ws.Range("E3").Value = "  +0.46%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'507.61"
$ws.Range("E5").Value = "  +0.66%  "
$ws.Range("D6").Value = "'153.41"
$ws.Range("E6").Value = "  -1.54%  "
$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("E8").Value = "  -6.05%  "
$ws.Range("D9").Value = "2.590.58"
$ws.Range("E9").Value = "  +0.73%  "
$ws.Range("D10").Value = "'6.58"
$ws.Range("E10").Value = "  +7.49%  "
$ws.Range("E11").Value = "  +1.14%  "
$ws.Range("E12").Value = "  +2.17%  "
$ws.Range("E13").Value = "  +0.79%  "
$ws.Range("D14").Value = "3.040.23"
$ws.Range("E14").Value = "  +0.65%  "
$ws.Range("D15").Value = "60.482.64"
$ws.Range("E15").Value = "  +1.84%  "
$ws.Range("D16").Value = "'21.58"
$ws.Range("E16").Value = "  -0.98%  "
$ws.Range("E17").Value = "  +2.56%  "
$ws.Range("D18").Value = "2.593.34"
$ws.Range("E18").Value = "  +0.62%  "
$ws.Range("E19").Value = "  +0.56%  "
$ws.Range("D20").Value = "'345.49"
$ws.Range("E20").Value = "  +2.55%  "
$ws.Range("E21").Value = "  +0.92%  "
$ws.Range("E22").Value = "  +1.65%  "
$ws.Range("E23").Value = "  -0.41%  "
$ws.Range("D24").Value = "'59.91"
$ws.Range("E24").Value = "  -0.07%  "
$ws.Range("E25").Value = "  +1.03%  "
$ws.Range("E26").Value = "  +0.95%  "
$ws.Range("E27").Value = "  -0.32%  "
$ws.Range("D28").Value = "0.0₃0844"
$ws.Range("E28").Value = "  +0.89%  "
$ws.Range("E29").Value = "  +0.11%  "
$ws.Range("E30").Value = "  -0.06%  "
$ws.Range("E31").Value = "  -0.25%  "
$ws.Range("D32").Value = "'153.79"
$ws.Range("E32").Value = "  -2.13%  "
$ws.Range("E33").Value = "  -0.61%  "
$ws.Range("E34").Value = "  +3.66%  "
$ws.Range("E35").Value = "  +1.76%  "
$ws.Range("E36").Value = "  -0.55%  "
$ws.Range("D37").Value = "'0.851"
$ws.Range("E37").Value = "  +8.58%  "
$ws.Range("D38").Value = "'0.851"
$ws.Range("E38").Value = "  +0.05%  "
$ws.Range("D39").Value = "'1.47"
$ws.Range("E39").Value = "  +2.50%  "
$ws.Range("E40").Value = "  +0.12%  "
$ws.Range("E41").Value = "  +2.23%  "
$ws.Range("D42").Value = "'294.52"
$ws.Range("E42").Value = "  -0.61%  "
$ws.Range("E43").Value = "  -0.32%  "
$ws.Range("D44").Value = "'0.0994"
$ws.Range("E44").Value = "  -2.20%  "
$ws.Range("B45").Value = "Hedera"
$ws.Range("C45").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D45").Value = "'0.0557"
$ws.Range("E45").Value = "  -2.89%  "
$ws.Range("B46").Value = "FirstDigitalUSD"
$ws.Range("C46").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D46").Value = "'0.996"
$ws.Range("E46").Value = "  -0.32%  "
$ws.Range("E47").Value = "  +3.89%  "
$ws.Range("D48").Value = "'4.85"
$ws.Range("E48").Value = "  -0.70%  "
$ws.Range("E49").Value = "  -1.21%  "
$ws.Range("D50").Value = "'10.29"
$ws.Range("E50").Value = "  +0.31%  "
$ws.Range("D51").Value = "2.002.05"
$ws.Range("E51").Value = "  +0.63%  "
